{"js": "// Remove the redundant \"then \" from the \"Repay a loan\" instructions so\n// that the sentence reads \"...select a loan, select one of your account,\n// then click repay loan button.\" instead of \"...select a loan, then\n// select one of your account, then click repay loan button.\"\nconst body = context.document.body;\n\nconst target = \"select a loan, then select one of your account\";\nconst replacement = \"select a loan, select one of your account\";\n\nconst results = body.search(target, { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text not found: \" + target);\n}\n\n// There is exactly one occurrence of this phrase in the document; replace it.\nconst found = results.items[0];\nfound.insertText(replacement, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Remove the redundant \"then \" from the \"Repay a loan\" instructions so\n# that the sentence reads \"...select a loan, select one of your account,\n# then click repay loan button.\" instead of \"...select a loan, then\n# select one of your account, then click repay loan button.\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"select a loan, then select one of your account\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"select a loan, select one of your account\"\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n# MatchCase=True, MatchWholeWord=False, MatchWildcards=False, MatchSoundsLike=False,\n# MatchAllWordForms=False, Forward=True, Wrap=wdFindContinue, Format=False,\n# Replace=wdReplaceOne (replace the single occurrence).\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceOne) | Out-Null\n"}
